# Update the single data row (row 2) of the 300922 profit-statement sheet
# with the latest quarterly figures (bot data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- REPORT_TYPE_CODE (J2): stays a text code, must remain "003" (not numeric 3) ---
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "003"

# --- REPORT_DATE (N2): quarter-end date, stored as plain text like the source data ---
$ws.Range("N2").Value = "2020-03-31 00:00:00"

# --- Numeric financial figures (row 2) ---
$ws.Range("O2").Value = 7534365.99
$ws.Range("P2").Value = 29278553.71
$ws.Range("Q2").Value = 23643696.46
$ws.Range("R2").Value = -28.6593858958
$ws.Range("S2").Value = 19448142.08
$ws.Range("T2").Value = 19448142.08
$ws.Range("U2").Value = -23.0346406466
$ws.Range("V2").Value = 474404.4
$ws.Range("W2").Value = 2040369.19
$ws.Range("X2").Value = -455105.8
$ws.Range("Y2").Value = 8866361.75
$ws.Range("Z2").Value = 8866417.93
$ws.Range("AA2").Value = 1332051.94
$ws.Range("AG2").Value = 245627.13
$ws.Range("AP2").Value = -17.5721818637
$ws.Range("AQ2").Value = 201.772376128175
$ws.Range("AR2").Value = 201.690629158218
$ws.Range("AS2").Value = 6185665.99
$ws.Range("AT2").Value = 205.989013195771
